$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "Bitcoin price change is fueled by global monetary growth (with forecast).xlsx",
  "Bank credit and M2 U.S.xlsx",
  "U.S GDP and GDI.xlsx",
  "US financial conditions, employment & equities.xlsx",
  "Real Incomes U.S (BEA).xlsx",
  "Bitty, SPX & GM2 fitted trends.xlsx",
  "Monetary Expansion and Risk Asset Activity.xlsx",
  "US Equity Indices.xlsx",
  "US and Global Monetary Aggregates.xlsx",
  "Financial conditions & employment, US.xlsx",
  "Bitcoin price change is fueled by global monetary growth.xlsx",
  "Personal savings U.S (BEA).xlsx",
  "Monetary Aggregates & Inflation, USA.xlsx",
  "Total construction spending in the U.S, nominal and real.xlsx",
  "US Gross Domestic Income and Equity Indices.xlsx"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
